$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 69.42974966666667
$ws.Range("H2").Value = 208.289249
$ws.Range("I2").Value = 0.6762611189535584
$ws.Range("J2").Value = 0.6762611189535584
$ws.Range("M2").Value = 101.8783343333333
$ws.Range("N2").Value = 305.635003
$ws.Range("O2").Value = 0.9137375742483709
$ws.Range("P2").Value = 0.913737574248371
$ws.Range("Q2").Value = 7073.387249220305
$ws.Range("R2").Value = 63660.48524298274
$ws.Range("S2").Value = 0.6179251943911135
$ws.Range("T2").Value = 0.6179251943911136

# Row 3
$ws.Range("G3").Value = 69.42974966666667
$ws.Range("H3").Value = 208.289249
$ws.Range("I3").Value = 0.6762611189535584
$ws.Range("J3").Value = 0.6762611189535584
$ws.Range("M3").Value = 0.050239
$ws.Range("O3").Value = 0.0004505890510780002
$ws.Range("P3").Value = 0.0004505890510780002
$ws.Range("Q3").Value = 3.488081193503667
$ws.Range("R3").Value = 31.392730741533
$ws.Range("S3").Value = 0.0003047158558702305
$ws.Range("T3").Value = 0.0003047158558702305

# Row 4
$ws.Range("G4").Value = 69.42974966666667
$ws.Range("H4").Value = 208.289249
$ws.Range("I4").Value = 0.6762611189535584
$ws.Range("J4").Value = 0.6762611189535584
$ws.Range("M4").Value = 3.444232
$ws.Range("N4").Value = 10.332696
$ws.Range("O4").Value = 0.03089100556484967
$ws.Range("P4").Value = 0.03089100556484967
$ws.Range("Q4").Value = 239.1321655539227
$ws.Range("R4").Value = 2152.189489985304
$ws.Range("S4").Value = 0.02089038598888584
$ws.Range("T4").Value = 0.02089038598888584

# Row 5
$ws.Range("G5").Value = 69.42974966666667
$ws.Range("H5").Value = 208.289249
$ws.Range("I5").Value = 0.6762611189535584
$ws.Range("J5").Value = 0.6762611189535584
$ws.Range("M5").Value = 0.07069366666666667
$ws.Range("N5").Value = 0.212081
$ws.Range("O5").Value = 0.0006340451079949399
$ws.Range("P5").Value = 0.00063404510799494
$ws.Range("Q5").Value = 4.908243579685444
$ws.Range("R5").Value = 44.17419221716901
$ws.Range("S5").Value = 0.0004287800541996878
$ws.Range("T5").Value = 0.0004287800541996879

# Row 6
$ws.Range("G6").Value = 69.42974966666667
$ws.Range("H6").Value = 208.289249
$ws.Range("I6").Value = 0.6762611189535584
$ws.Range("J6").Value = 0.6762611189535584
$ws.Range("M6").Value = 1.402613333333333
$ws.Range("N6").Value = 4.20784
$ws.Range("O6").Value = 0.01257991223742545
$ws.Range("P6").Value = 0.01257991223742545
$ws.Range("Q6").Value = 97.38309261246222
$ws.Range("R6").Value = 876.4478335121601
$ws.Range("S6").Value = 0.0085073055260189
$ws.Range("T6").Value = 0.008507305526018902

# Row 7
$ws.Range("G7").Value = 69.42974966666667
$ws.Range("H7").Value = 208.289249
$ws.Range("I7").Value = 0.6762611189535584
$ws.Range("J7").Value = 0.6762611189535584
$ws.Range("M7").Value = 4.650161000000001
$ws.Range("N7").Value = 13.950483
$ws.Range("O7").Value = 0.04170687379028095
$ws.Range("P7").Value = 0.04170687379028095
$ws.Range("Q7").Value = 322.8595141396964
$ws.Range("R7").Value = 2905.735627257267
$ws.Range("S7").Value = 0.02820473713747023
$ws.Range("T7").Value = 0.02820473713747023

# Row 8
$ws.Range("I8").Value = 0.1530049204123442
$ws.Range("J8").Value = 0.1530049204123442
$ws.Range("M8").Value = 101.8783343333333
$ws.Range("N8").Value = 305.635003
$ws.Range("O8").Value = 0.9137375742483709
$ws.Range("P8").Value = 0.913737574248371
$ws.Range("Q8").Value = 1600.362674683012
$ws.Range("R8").Value = 14403.26407214711
$ws.Range("S8").Value = 0.1398063448256405
$ws.Range("T8").Value = 0.1398063448256405

# Row 9
$ws.Range("I9").Value = 0.1530049204123442
$ws.Range("J9").Value = 0.1530049204123442
$ws.Range("M9").Value = 0.050239
$ws.Range("O9").Value = 0.0004505890510780002
$ws.Range("P9").Value = 0.0004505890510780002
$ws.Range("S9").Value = 0.00006894234189886312
$ws.Range("T9").Value = 0.00006894234189886312

# Row 10
$ws.Range("I10").Value = 0.1530049204123442
$ws.Range("J10").Value = 0.1530049204123442
$ws.Range("M10").Value = 3.444232
$ws.Range("N10").Value = 10.332696
$ws.Range("O10").Value = 0.03089100556484967
$ws.Range("P10").Value = 0.03089100556484967
$ws.Range("Q10").Value = 54.10395028362134
$ws.Range("R10").Value = 486.9355525525921
$ws.Range("S10").Value = 0.004726475847907106
$ws.Range("T10").Value = 0.004726475847907107

# Row 11
$ws.Range("I11").Value = 0.1530049204123442
$ws.Range("J11").Value = 0.1530049204123442
$ws.Range("M11").Value = 0.07069366666666667
$ws.Range("N11").Value = 0.212081
$ws.Range("O11").Value = 0.0006340451079949399
$ws.Range("P11").Value = 0.00063404510799494
$ws.Range("Q11").Value = 1.110496222873556
$ws.Range("R11").Value = 9.994466005862002
$ws.Range("S11").Value = 0.00009701202128660198
$ws.Range("T11").Value = 0.000097012021286602

# Row 12
$ws.Range("I12").Value = 0.1530049204123442
$ws.Range("J12").Value = 0.1530049204123442
$ws.Range("M12").Value = 1.402613333333333
$ws.Range("N12").Value = 4.20784
$ws.Range("O12").Value = 0.01257991223742545
$ws.Range("P12").Value = 0.01257991223742545
$ws.Range("Q12").Value = 22.03304598929778
$ws.Range("R12").Value = 198.29741390368
$ws.Range("S12").Value = 0.001924788470681557
$ws.Range("T12").Value = 0.001924788470681557

# Row 13
$ws.Range("I13").Value = 0.1530049204123442
$ws.Range("J13").Value = 0.1530049204123442
$ws.Range("M13").Value = 4.650161000000001
$ws.Range("N13").Value = 13.950483
$ws.Range("O13").Value = 0.04170687379028095
$ws.Range("P13").Value = 0.04170687379028095
$ws.Range("Q13").Value = 73.04736717934068
$ws.Range("R13").Value = 657.4263046140661
$ws.Range("S13").Value = 0.006381356904929622
$ws.Range("T13").Value = 0.006381356904929622

# Row 14
$ws.Range("G14").Value = 17.31398133333333
$ws.Range("H14").Value = 51.941944
$ws.Range("I14").Value = 0.1686420078746507
$ws.Range("J14").Value = 0.1686420078746507
$ws.Range("M14").Value = 101.8783343333333
$ws.Range("N14").Value = 305.635003
$ws.Range("O14").Value = 0.9137375742483709
$ws.Range("P14").Value = 0.913737574248371
$ws.Range("Q14").Value = 1763.919578918426
$ws.Range("R14").Value = 15875.27621026583
$ws.Range("S14").Value = 0.154094539191758
$ws.Range("T14").Value = 0.154094539191758

# Row 15
$ws.Range("G15").Value = 17.31398133333333
$ws.Range("H15").Value = 51.941944
$ws.Range("I15").Value = 0.1686420078746507
$ws.Range("J15").Value = 0.1686420078746507
$ws.Range("M15").Value = 0.050239
$ws.Range("O15").Value = 0.0004505890510780002
$ws.Range("P15").Value = 0.0004505890510780002
$ws.Range("Q15").Value = 0.8698371082053333
$ws.Range("R15").Value = 7.828533973848
$ws.Range("S15").Value = 0.00007598824230012748
$ws.Range("T15").Value = 0.00007598824230012748

# Row 16
$ws.Range("G16").Value = 17.31398133333333
$ws.Range("H16").Value = 51.941944
$ws.Range("I16").Value = 0.1686420078746507
$ws.Range("J16").Value = 0.1686420078746507
$ws.Range("M16").Value = 3.444232
$ws.Range("N16").Value = 10.332696
$ws.Range("O16").Value = 0.03089100556484967
$ws.Range("P16").Value = 0.03089100556484967
$ws.Range("Q16").Value = 59.63336855566934
$ws.Range("R16").Value = 536.700317001024
$ws.Range("S16").Value = 0.005209521203723255
$ws.Range("T16").Value = 0.005209521203723256

# Row 17
$ws.Range("G17").Value = 17.31398133333333
$ws.Range("H17").Value = 51.941944
$ws.Range("I17").Value = 0.1686420078746507
$ws.Range("J17").Value = 0.1686420078746507
$ws.Range("M17").Value = 0.07069366666666667
$ws.Range("N17").Value = 0.212081
$ws.Range("O17").Value = 0.0006340451079949399
$ws.Range("P17").Value = 0.00063404510799494
$ws.Range("Q17").Value = 1.223988825051556
$ws.Range("R17").Value = 11.015899425464
$ws.Range("S17").Value = 0.0001069266400953664
$ws.Range("T17").Value = 0.0001069266400953664

# Row 18
$ws.Range("G18").Value = 17.31398133333333
$ws.Range("H18").Value = 51.941944
$ws.Range("I18").Value = 0.1686420078746507
$ws.Range("J18").Value = 0.1686420078746507
$ws.Range("M18").Value = 1.402613333333333
$ws.Range("N18").Value = 4.20784
$ws.Range("O18").Value = 0.01257991223742545
$ws.Range("P18").Value = 0.01257991223742545
$ws.Range("Q18").Value = 24.28482107121778
$ws.Range("R18").Value = 218.56338964096
$ws.Range("S18").Value = 0.002121501658606318
$ws.Range("T18").Value = 0.002121501658606318

# Row 19
$ws.Range("G19").Value = 17.31398133333333
$ws.Range("H19").Value = 51.941944
$ws.Range("I19").Value = 0.1686420078746507
$ws.Range("J19").Value = 0.1686420078746507
$ws.Range("M19").Value = 4.650161000000001
$ws.Range("N19").Value = 13.950483
$ws.Range("O19").Value = 0.04170687379028095
$ws.Range("P19").Value = 0.04170687379028095
$ws.Range("Q19").Value = 80.51280075099469
$ws.Range("R19").Value = 724.6152067589521
$ws.Range("S19").Value = 0.007033530938167621
$ws.Range("T19").Value = 0.007033530938167621

# Row 20
$ws.Range("G20").Value = 0.2147746666666667
$ws.Range("H20").Value = 0.644324
$ws.Range("I20").Value = 0.002091952759446708
$ws.Range("J20").Value = 0.002091952759446708
$ws.Range("M20").Value = 101.8783343333333
$ws.Range("N20").Value = 305.635003
$ws.Range("O20").Value = 0.9137375742483709
$ws.Range("P20").Value = 0.913737574248371
$ws.Range("Q20").Value = 21.88088529699689
$ws.Range("R20").Value = 196.927967672972
$ws.Range("S20").Value = 0.001911495839859021
$ws.Range("T20").Value = 0.001911495839859021

# Row 21
$ws.Range("G21").Value = 0.2147746666666667
$ws.Range("H21").Value = 0.644324
$ws.Range("I21").Value = 0.002091952759446708
$ws.Range("J21").Value = 0.002091952759446708
$ws.Range("M21").Value = 0.050239
$ws.Range("O21").Value = 0.0004505890510780002
$ws.Range("P21").Value = 0.0004505890510780002
$ws.Range("Q21").Value = 0.01079006447866667
$ws.Range("R21").Value = 0.09711058030799999
$ws.Range("S21").Value = 0.0000009426110087790964
$ws.Range("T21").Value = 0.0000009426110087790964

# Row 22
$ws.Range("G22").Value = 0.2147746666666667
$ws.Range("H22").Value = 0.644324
$ws.Range("I22").Value = 0.002091952759446708
$ws.Range("J22").Value = 0.002091952759446708
$ws.Range("M22").Value = 3.444232
$ws.Range("N22").Value = 10.332696
$ws.Range("O22").Value = 0.03089100556484967
$ws.Range("P22").Value = 0.03089100556484967
$ws.Range("Q22").Value = 0.7397337797226666
$ws.Range("R22").Value = 6.657604017504
$ws.Range("S22").Value = 0.0000646225243334709
$ws.Range("T22").Value = 0.0000646225243334709

# Row 23
$ws.Range("G23").Value = 0.2147746666666667
$ws.Range("H23").Value = 0.644324
$ws.Range("I23").Value = 0.002091952759446708
$ws.Range("J23").Value = 0.002091952759446708
$ws.Range("M23").Value = 0.07069366666666667
$ws.Range("N23").Value = 0.212081
$ws.Range("O23").Value = 0.0006340451079949399
$ws.Range("P23").Value = 0.00063404510799494
$ws.Range("Q23").Value = 0.01518320869377778
$ws.Range("R23").Value = 0.136648878244
$ws.Range("S23").Value = 0.000001326392413283701
$ws.Range("T23").Value = 0.000001326392413283701

# Row 24
$ws.Range("G24").Value = 0.2147746666666667
$ws.Range("H24").Value = 0.644324
$ws.Range("I24").Value = 0.002091952759446708
$ws.Range("J24").Value = 0.002091952759446708
$ws.Range("M24").Value = 1.402613333333333
$ws.Range("N24").Value = 4.20784
$ws.Range("O24").Value = 0.01257991223742545
$ws.Range("P24").Value = 0.01257991223742545
$ws.Range("Q24").Value = 0.3012458111288889
$ws.Range("R24").Value = 2.71121230016
$ws.Range("S24").Value = 0.00002631658211867959
$ws.Range("T24").Value = 0.00002631658211867959

# Row 25
$ws.Range("G25").Value = 0.2147746666666667
$ws.Range("H25").Value = 0.644324
$ws.Range("I25").Value = 0.002091952759446708
$ws.Range("J25").Value = 0.002091952759446708
$ws.Range("M25").Value = 4.650161000000001
$ws.Range("N25").Value = 13.950483
$ws.Range("O25").Value = 0.04170687379028095
$ws.Range("P25").Value = 0.04170687379028095
$ws.Range("Q25").Value = 0.9987367787213335
$ws.Range("R25").Value = 8.988631008492002
$ws.Range("S25").Value = 0.00008724880971347383
$ws.Range("T25").Value = 0.00008724880971347383
